$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect new "through" date
$ws.Name = "Through 2021-12-18"

# Update the label in column A for the December row
$ws.Range("A13").Value = "December (through 12-18)"

# Update December figures (row 13) for each year column B..H
$ws.Range("B13").Value = 23
$ws.Range("C13").Value = 57
$ws.Range("D13").Value = 70
$ws.Range("E13").Value = 42
$ws.Range("F13").Value = 29
$ws.Range("G13").Value = 83
$ws.Range("H13").Value = 130

# Update Total figures (row 14) for each year column B..H
$ws.Range("B14").Value = 314
$ws.Range("C14").Value = 620
$ws.Range("D14").Value = 891
$ws.Range("E14").Value = 724
$ws.Range("F14").Value = 563
$ws.Range("G14").Value = 1347
$ws.Range("H14").Value = 1773
